$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 21.827068328857422
$ws.Range("C2").Value = 5.896552085876465
$ws.Range("D2").Value = 15.231202125549316
$ws.Range("E2").Value = 57.85714340209961
